# Auto-generated edit script: updates cryptos Price (D) and Volume(1h) (E) columns
# to match the new scraped values from the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, [string]$Text)
    # Force the cell to be treated as text so numeric-looking strings
    # (e.g. "196.17") are not coerced into Number cells, then strip the
    # temporary text format so no stray style is left behind.
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.ClearFormats()
}

Set-TextValue $ws.Range("D2") "69.598.80"
Set-TextValue $ws.Range("E2") "  -0.92%  "
Set-TextValue $ws.Range("D3") "3.549.35"
Set-TextValue $ws.Range("E3") "  -1.42%  "
Set-TextValue $ws.Range("E4") "  +0.13%  "
Set-TextValue $ws.Range("D5") "196.17"
Set-TextValue $ws.Range("E5") "  +0.10%  "
Set-TextValue $ws.Range("D6") "585.82"
Set-TextValue $ws.Range("E7") "  -2.23%  "
Set-TextValue $ws.Range("E8") "  +0.04%  "
Set-TextValue $ws.Range("D9") "0.206"
Set-TextValue $ws.Range("E9") "  -0.08%  "
Set-TextValue $ws.Range("E10") "  -2.54%  "
Set-TextValue $ws.Range("D11") "52.86"
Set-TextValue $ws.Range("E11") "  -1.72%  "
Set-TextValue $ws.Range("E12") "  -5.07%  "
Set-TextValue $ws.Range("D13") "9.24"
Set-TextValue $ws.Range("E13") "  -3.24%  "
Set-TextValue $ws.Range("D14") "4.116.90"
Set-TextValue $ws.Range("E14") "  -1.52%  "
Set-TextValue $ws.Range("D15") "666.38"
Set-TextValue $ws.Range("E15") "  +11.81%  "
Set-TextValue $ws.Range("D16") "69.735.22"
Set-TextValue $ws.Range("E16") "  -0.86%  "
Set-TextValue $ws.Range("D17") "3.567.43"
Set-TextValue $ws.Range("E17") "  -1.34%  "
Set-TextValue $ws.Range("D18") "12.50"
Set-TextValue $ws.Range("E18") "  -4.45%  "
Set-TextValue $ws.Range("E19") "  -0.74%  "
Set-TextValue $ws.Range("D20") "18.44"
Set-TextValue $ws.Range("E20") "  -3.21%  "
Set-TextValue $ws.Range("E21") "  -2.99%  "
Set-TextValue $ws.Range("D22") "18.03"
Set-TextValue $ws.Range("E22") "  +1.10%  "
Set-TextValue $ws.Range("E23") "  +3.72%  "
Set-TextValue $ws.Range("D24") "105.30"
Set-TextValue $ws.Range("E24") "  +3.18%  "
Set-TextValue $ws.Range("E25") "  -4.85%  "
Set-TextValue $ws.Range("D26") "2.93"
Set-TextValue $ws.Range("E26") "  -3.34%  "
Set-TextValue $ws.Range("D27") "10.15"
Set-TextValue $ws.Range("E27") "  -5.42%  "
Set-TextValue $ws.Range("D28") "9.59"
Set-TextValue $ws.Range("E28") "  -0.46%  "
Set-TextValue $ws.Range("D29") "33.39"
Set-TextValue $ws.Range("E29") "  -1.15%  "
Set-TextValue $ws.Range("D30") "4.36"
Set-TextValue $ws.Range("E30") "  -8.99%  "
Set-TextValue $ws.Range("D31") "6.78"
Set-TextValue $ws.Range("E31") "  -5.31%  "
Set-TextValue $ws.Range("D32") "11.78"
Set-TextValue $ws.Range("E32") "  -4.09%  "
Set-TextValue $ws.Range("E33") "  -4.50%  "
Set-TextValue $ws.Range("D34") "62.03"
Set-TextValue $ws.Range("E34") "  -1.85%  "
Set-TextValue $ws.Range("D35") "3.795.11"
Set-TextValue $ws.Range("E35") "  -2.68%  "
Set-TextValue $ws.Range("D36") "3.76"
Set-TextValue $ws.Range("E36") "  +6.75%  "
Set-TextValue $ws.Range("D37") "0.0₃0809"
Set-TextValue $ws.Range("E37") "  -10.45%  "
Set-TextValue $ws.Range("E38") "  +0.10%  "
Set-TextValue $ws.Range("D39") "499.06"
Set-TextValue $ws.Range("E39") "  -4.43%  "
Set-TextValue $ws.Range("D40") "2.89"
Set-TextValue $ws.Range("E40") "  -6.57%  "
Set-TextValue $ws.Range("E41") "  -4.68%  "
Set-TextValue $ws.Range("D42") "0.134"
Set-TextValue $ws.Range("E42") "  +0.98%  "
Set-TextValue $ws.Range("E43") "  -6.20%  "
Set-TextValue $ws.Range("D44") "0.0451"
Set-TextValue $ws.Range("E44") "  -0.56%  "
Set-TextValue $ws.Range("D45") "2.90"
Set-TextValue $ws.Range("E45") "  +1.61%  "
Set-TextValue $ws.Range("E46") "  -1.25%  "
Set-TextValue $ws.Range("E47") "  -2.42%  "
Set-TextValue $ws.Range("E48") "  -0.27%  "
Set-TextValue $ws.Range("E49") "  -2.94%  "
Set-TextValue $ws.Range("E50") "  +19.67%  "
Set-TextValue $ws.Range("D51") "2.69"
Set-TextValue $ws.Range("E51") "  +61.62%  "
